$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.230.07'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +2.23%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.097.94'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +4.11%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '251.13'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.31%  '

$ws.Range("E6").Value = '  +0.47%  '

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '54.09'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +20.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '61.52'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.55%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.374'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.68%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0744'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +4.42%  '

$ws.Range("E12").Value = '  +7.29%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.34'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +5.44%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.403.67'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +4.00%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.841'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +3.86%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.107.73'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +4.59%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.18'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +5.85%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.161.85'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +2.34%  '

$ws.Range("E19").Value = '  +1.80%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.56'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +12.95%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0838'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.77%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '241.89'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +2.48%  '

$ws.Range("E23").Value = '  +7.69%  '

$ws.Range("E24").Value = '  -0.05%  '

$ws.Range("E25").Value = '  +1.25%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '170.88'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +4.42%  '

$ws.Range("E27").Value = '  +7.87%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.68'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +5.50%  '

$ws.Range("E29").Value = '  +4.10%  '

$ws.Range("E30").Value = '  +1.13%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.08'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +27.46%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.53'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.81%  '

$ws.Range("B33").Value = 'Gas'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.63'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -4.12%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0614'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +4.34%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0905'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +12.45%  '

$ws.Range("E36").Value = '  -0.08%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.31'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +7.18%  '

$ws.Range("E38").Value = '  -0.93%  '

$ws.Range("E39").Value = '  +3.01%  '

$ws.Range("E40").Value = '  +1.18%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.32'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +14.88%  '

$ws.Range("E42").Value = '  +4.61%  '

$ws.Range("E43").Value = '  +4.87%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.01'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +3.22%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0915'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +12.11%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.78'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.04%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.321.58'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.38%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.97'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +7.14%  '

$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.06'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +14.69%  '

$ws.Range("B50").Value = 'FTXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.83'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +85.36%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.290.64'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +3.84%  '
